$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Double the values in A2:A26
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    $cell.Value = $v * 2
}

# Update the view: scroll so A7 is the top-left cell, and select G19
$ws.Range("G19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
